# Cost updates on the P2G sheet: column C (Variable O&M) switches from a flat
# numeric 3700 to a declining series of values that are stored as *text*
# (shared-string) cells, matching a paste-as-text / CSV-import style edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("P2G")

# Row 2 (year 2019) through row 33 (year 2050) - the new text values, in order.
$values = @(
    "6809","6614","6420","6225","6031","5836","5642","5447",
    "5253","5058","4864","4669","4475","4280","4086","3891",
    "3696","3696","3696","3696","3696","3696","3696","3696",
    "3696","3696","3696","3696","3696","3696","3696","3696"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 3)

    # Writing the digit-string straight into .Value/.Formula gets reinterpreted
    # as a number by Excel. Route it through a text-producing formula and then
    # paste-special just the value back in, which is how Excel preserves a
    # genuine text cell (shared string) without touching any cell styles.
    $cell.Formula = '="' + $values[$i] + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Restore the active selection like the source workbook shows.
$ws.Range("H5").Select()
